# edit.ps1 -- reproduce the tcn_p115v.docx commit ("fixes per #20; regenerated files")
#
# 1. Split "Ne retire pas" so its leading "N" carries a new comment
#    (author "Soersha Dyon") referencing the recipe "Recuire se moules
#    à noyau".
# 2. Drop the comma in "soient essorés et affermis, car aultrement".
# 3. Turn the semicolon after "foeu" into a comma.
# 4. Lower-case the "I" in ", Ilz se recuisent bien mieulx dans un".
# 5. Drop the trailing comma in "recuire deulx foys,".

$d = $word.ActiveDocument

# The comment must be authored by "Soersha Dyon" -- Word stamps new
# comments with Application.UserName, so set it before calling Add().
$word.UserName = "Soersha Dyon"

# --- 1. comment on the leading "N" of "Ne retire pas" -----------------
$rng = $d.Content
$rng.Find.Execute("Ne retire pas", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$anchor = $d.Range($rng.Start, $rng.Start + 1)
$d.Comments.Add($anchor, "For markup : this is part of the recipe ""Recuire se moules à noyau""") | Out-Null

# --- 2. remove the comma before "car aultrement" -----------------------
$d.Content.Find.Execute("soient essorés et affermis, car aultrement", `
                         $true, $false, $false, $false, $false, $true, 1, `
                         $false, "soient essorés et affermis car aultrement", 2) | Out-Null

# --- 3. "foeu;" -> "foeu," ----------------------------------------------
$d.Content.Find.Execute("ce mesme foeu;", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "ce mesme foeu,", 2) | Out-Null

# --- 4. "Ilz" -> "ilz" ---------------------------------------------------
$d.Content.Find.Execute(", Ilz se recuisent bien mieulx dans un", $true, `
                         $false, $false, $false, $false, $true, 1, $false, `
                         ", ilz se recuisent bien mieulx dans un", 2) | Out-Null

# --- 5. drop the trailing comma in "recuire deulx foys," ----------------
$d.Content.Find.Execute("recuire deulx foys,", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "recuire deulx foys", 2) | Out-Null
